$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values (duplicates removed from results)
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 21
$ws.Range("I2").Value = 40
$ws.Range("J2").Value = 13
$ws.Range("E3").Value = 101
$ws.Range("H3").Value = 101
$ws.Range("I3").Value = 53
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 6
$ws.Range("I4").Value = 20
$ws.Range("J4").Value = 9
$ws.Range("E5").Value = 3
$ws.Range("H5").Value = 7
$ws.Range("I5").Value = 9
$ws.Range("C6").Value = 5
$ws.Range("E6").Value = 18
$ws.Range("I6").Value = 658
$ws.Range("F7").Value = 1
$ws.Range("I7").Value = 42
$ws.Range("J7").Value = 22
$ws.Range("I8").Value = 9
$ws.Range("I9").Value = 10
$ws.Range("J9").Value = 1

# Freeze the header row (row 1) and set the active selection
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B18").Select()
